$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new user record (P1004 / Password) as row 10, matching columns
# "Hospital ID" (A) and "Password" (B) used by the rest of the sheet.
$ws.Range("A10").Value = "P1004"
$ws.Range("B10").Value = "Password"
